$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "Team member" header (column E) and add the new
# "Team member Actual" header/data in column F.
$ws.Range("E3").Value = "Team member Initial"

$ws.Range("F3").Value  = "Team member Actual"
$ws.Range("F5").Value  = "Tony Huynh"
$ws.Range("F6").Value  = "Andrew Lam"
$ws.Range("F7").Value  = "Tony Huynh"
$ws.Range("F8").Value  = "Tony Huynh"
$ws.Range("F9").Value  = "Tony Huynh/Andrew Lam"
$ws.Range("F10").Value = "Tony Huynh/Andrew Lam"
$ws.Range("F11").Value = "Tony Huynh/Andrew Lam"
$ws.Range("F12").Value = "Andrew Lam"
$ws.Range("F13").Value = "Tony Huynh"
$ws.Range("F14").Value = "Tony Huynh"
$ws.Range("F15").Value = "Andrew Lam"
$ws.Range("F16").Value = "Tony Huynh"
$ws.Range("F17").Value = "Nathaniel Leake"
$ws.Range("F22").Value = "Tony Huynh/Andrew Lam"
$ws.Range("F25").Value = "Andrew Lam/Tony Huyhn"
$ws.Range("F31").Value = "Tony Huynh"

$ws.Range("F15").Select()
